$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("test")

# Clear out the stale platform/device data cells in rows 3-12
$clearCells = @(
    "B3","C3","D3","E3","F3",
    "B4","C4","D4","E4","F4",
    "B5","C5","D5","E5","F5",
    "B6","C6","D6","E6","F6",
    "B7","F7",
    "B8","C8","F8",
    "B9","C9","F9",
    "C10","F10",
    "B11","C11","F11",
    "B12","C12","F12","M12"
)
foreach ($addr in $clearCells) {
    $ws2.Range($addr).Value = ""
}

# Add new appPackage/appActivity/bundleId style row for BrowserStack (Google Drive app)
$ws2.Range("K15").Value = "com.google.android.apps.docs"
$ws2.Range("L15").Value = ".drive.startup.StartupActivity"
$ws2.Range("M15").Value = "com.apple.preferences"
$ws2.Range("M15").NumberFormat = "@"

# Update the view: selection now on L15, no frozen/topLeft override
$ws2.Activate() | Out-Null
$ws2.Range("L15").Select() | Out-Null
